$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kanjiSentence = "彼は今、新薬の研究開発に|挑んで|いる。"
$sosyouSentence = "住民が建設会社を相手に、|訴訟|を起こした。"
$kanjiSentenceLong = "彼は今、新薬の研究開発に|挑んで|いるうううううううううううううううううううううううううううううううううううううううう。"

# Every cell that previously held '彼は今、新薬の研究開発に・挑んで・いる。'
# now holds the same sentence with '|' instead of '・', except B3 which is
# replaced by a new, much longer variant of the sentence.
$ws.Range("B3").Value = $kanjiSentenceLong
$ws.Range("B5").Value = $kanjiSentence
$ws.Range("B7").Value = $kanjiSentence
$ws.Range("B9").Value = $kanjiSentence
$ws.Range("B11").Value = $kanjiSentence
$ws.Range("B13").Value = $kanjiSentence
$ws.Range("B15").Value = $kanjiSentence
$ws.Range("B18").Value = $kanjiSentence
$ws.Range("B20").Value = $kanjiSentence
$ws.Range("B22").Value = $kanjiSentence

# Every cell that previously held '住民が建設会社を相手に、・訴訟・を起こした。'
# now holds the same sentence with '|' instead of '・'.
$ws.Range("B4").Value = $sosyouSentence
$ws.Range("B6").Value = $sosyouSentence
$ws.Range("B8").Value = $sosyouSentence
$ws.Range("B10").Value = $sosyouSentence
$ws.Range("B12").Value = $sosyouSentence
$ws.Range("B14").Value = $sosyouSentence
$ws.Range("B16").Value = $sosyouSentence
$ws.Range("B19").Value = $sosyouSentence
$ws.Range("B21").Value = $sosyouSentence
$ws.Range("B23").Value = $sosyouSentence

# Update the active cell selection shown in the sheet view.
$ws.Range("J8").Select()
